$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 74, shifting existing rows 74..194 down to 75..195
$ws.Rows.Item(74).Insert()

# Populate the new row 74 with the new data
$ws.Cells.Item(74, 1).Value = 10
$ws.Cells.Item(74, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(74, 3).Value = "La Araucanía"
$ws.Cells.Item(74, 4).Value = 44571
$ws.Cells.Item(74, 5).Value = 9
$ws.Cells.Item(74, 6).Value = 100112039
$ws.Cells.Item(74, 7).Value = "Ciboulette"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 65
$ws.Cells.Item(74, 11).Value = 5000
$ws.Cells.Item(74, 12).Value = 5000
$ws.Cells.Item(74, 13).Value = 5000
$ws.Cells.Item(74, 14).Value = "$/docena de atados"
$ws.Cells.Item(74, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(74, 16).Value = 1667
$ws.Cells.Item(74, 17).Value = 3
$ws.Cells.Item(74, 18).Value = "Hortaliza"
